# OLX Monitor 2026-02-22 12:16 refresh
#
# Appends a fresh monitoring snapshot (8 still-active listings) to the
# bottom of the "PODSUMOWANIE" sheet's listing log (rows 123-130), mirroring
# the previous snapshot block (rows 115-122) except for the "last checked"
# timestamp (column A) and one listing whose tracked counter (column D)
# reset to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "PODSUMOWANIE" is the first / active sheet

# 1) Clone the previous snapshot block (rows 115-122, cols A-H) -- values
#    first, then formats -- into the new block starting at row 123. This
#    reproduces both the per-listing data and the alternating cell styles
#    (s=13/14/15) without having to hardcode style ids.
$ws.Range("A115:H122").Copy()
$ws.Range("A123").PasteSpecial()

$ws.Range("A115:H122").Copy()
$ws.Range("A123").PasteSpecial(-4122)   # xlPasteFormats

# 2) Stamp the new block with this run's "last checked" timestamp.
$ws.Range("A123:A130").Value = "2026-02-22 12:16:18"

# 3) One listing's counter (column D) differs from its prior snapshot --
#    "WOLNY OD ZARAZ! Pokoj jedynka, ul. Romanowskiego 58" (row 127) is 0
#    this time instead of the 58640 it was cloned with.
$ws.Cells.Item(127, 4).Value = 0
